# Gantt chart update: replace placeholder activity rows 17-20 with updated
# task names/dates, and add 7 newly tracked tasks (rows 21-27), matching the
# new project plan ("Feature" coding tasks, updated reports/docs tasks).
# Rows 28-39 (generic "Activity NN" placeholders) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 17: #013 Writing the code for Feature 1 ---
$ws.Range("B17").Value = "#013     Writing the code for Feature 1 "
$ws.Range("C17").Value = 24
$ws.Range("D17").Value = 24
$ws.Range("E17").Value = 24
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 1

# --- Row 18: #014 Writing the code for Feature 2 ---
$ws.Range("B18").Value = "#014     Writing the code for Feature 2"
$ws.Range("C18").Value = 25
$ws.Range("D18").Value = 27
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 1

# --- Row 19: #015 Writing the code for Feature 3 ---
$ws.Range("B19").Value = "#015     Writing the code for Feature 3 "
$ws.Range("C19").Value = 27
$ws.Range("D19").Value = 29
$ws.Range("E19").Value = 27
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 1

# --- Row 20: #016 Writing the code for Feature 4 ---
$ws.Range("B20").Value = "#016     Writing the code for Feature 4 "
$ws.Range("C20").Value = 29
$ws.Range("D20").Value = 31
$ws.Range("E20").Value = 28
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 1

# --- Row 21: #017 Writing the code for Feature 5 ---
$ws.Range("B21").Value = "#017     Writing the code for Feature 5 "
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = 30
$ws.Range("F21").Value = 32
$ws.Range("G21").Value = 1

# --- Row 22: #018 Creating a Software Testing Report.docx ---
$ws.Range("B22").Value = "#018     Creating a Software Testing Report.docx "
$ws.Range("C22").Value = 33
$ws.Range("D22").Value = 36
$ws.Range("E22").Value = 32
$ws.Range("F22").Value = 36
$ws.Range("G22").Value = 1

# --- Row 23: #019 Preparing the Executive Summary.docx ---
$ws.Range("B23").Value = "#019     Preparing the Executive Summary.docx "
$ws.Range("C23").Value = 36
$ws.Range("D23").Value = 38
$ws.Range("E23").Value = 36
$ws.Range("F23").Value = 38
$ws.Range("G23").Value = 1

# --- Row 24: #020 Updating Project Plan.docx ---
$ws.Range("B24").Value = "#020     Updating Project Plan.docx "
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = 38
$ws.Range("F24").Value = 40
$ws.Range("G24").Value = 1

# --- Row 25: #021 Updating Software Design Document.docx ---
$ws.Range("B25").Value = "#021     Updating Software Design Document.docx "
$ws.Range("C25").Value = 39
$ws.Range("D25").Value = 40
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 1

# --- Row 26: #021 Updating Gantt Chart ---
$ws.Range("B26").Value = "#021     Updating Gantt Chart "
$ws.Range("C26").Value = 40
$ws.Range("D26").Value = 41
$ws.Range("E26").Value = 41
$ws.Range("F26").Value = 41
$ws.Range("G26").Value = 1

# --- Row 27: #022 Downloading git_log.txt ---
$ws.Range("B27").Value = "#022     Downloading git_log.txt "
$ws.Range("C27").Value = 41
$ws.Range("D27").Value = 41
$ws.Range("E27").Value = 41
$ws.Range("F27").Value = 41
$ws.Range("G27").Value = 1

# Update the view: scroll so row 19 is at the top and select G28 (matches the
# author's on-save cursor position), and reset zoom to 100%.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("G28").Select()
